# The source export's row for id 86420981 (Varglav / Letharia vulpina) and
# the row for id 110346530 (Skrovellav / Lobaria scrobiculata) traded
# places: everything that used to live in row 2 now lives in row 3, and
# vice versa. Re-create that by swapping the two rows' contents column by
# column, leaving the header row (row 1) alone.
#
# A straight bulk Range.Value2 swap works for almost every column, but a
# handful of text cells "look like" a number or an ISO date
# (e.g. "1", "2019-08-19") and Excel's usual typed-value coercion would
# silently turn them into a real number/date when written through
# Value2. Those few cells are therefore written back with a leading
# apostrophe, which is how Excel keeps text that resembles a number as
# literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Range($fromCol, $toCol) {
    $r2 = $ws.Range($fromCol + "2:" + $toCol + "2")
    $r3 = $ws.Range($fromCol + "3:" + $toCol + "3")
    $v2 = $r2.Value2
    $v3 = $r3.Value2
    $r2.Value2 = $v3
    $r3.Value2 = $v2
}

# Bulk-swap the columns that never risk an implicit number/date coercion.
# AT and AY are skipped on purpose: they are blank in both rows already,
# and round-tripping an already-blank cell through Value2 would needlessly
# touch it for no actual change.
Swap-Range "A" "H"
Swap-Range "J" "X"
Swap-Range "Z" "Z"
Swap-Range "AB" "AS"
Swap-Range "AU" "AX"

# Column I ("Antal"): row 2 had the literal text "1", row 3 was blank.
# After swapping, row 2 becomes blank and row 3 gets the text "1" - force
# the quoted literal so it is not re-typed as the number 1.
$ws.Range("I2").Value2 = ""
$ws.Range("I3").Value2 = "'1"

# Column Y ("Startdatum") and AA ("Slutdatum") hold ISO date strings as
# plain text; re-enter them with a leading apostrophe so they remain text
# instead of being parsed into date serials.
$ws.Range("Y2").Value2 = "'2023-06-25"
$ws.Range("AA2").Value2 = "'2023-06-25"
$ws.Range("Y3").Value2 = "'2019-08-19"
$ws.Range("AA3").Value2 = "'2019-08-19"
